# "updated printables for easier printing"
#
# 1) Bump the cached text of every "Update automatically" date field
#    (slideMaster, all 11 slideLayouts, notesMaster) from 3/3/24 to 11/13/24.
# 2) Reposition five of the six hexagon shapes on slide 1 so the grid
#    prints more cleanly.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders ("datetimeFigureOut" fields)
# ---------------------------------------------------------------------
$newDate = "11/13/24"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master
Update-DatePlaceholders $p.SlideMaster.Shapes

# Every slide layout
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholders $layouts.Item($L).Shapes
}

# Notes master (its Shapes collection doesn't accept direct text writes in
# this host, so go through the HeadersFooters facade instead)
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate

# ---------------------------------------------------------------------
# 2) Hexagon repositioning on slide 1
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

# Target positions in EMU -> points (Shape.Left / Shape.Top are in points).
# A couple of the EMU targets (x=3896227 / y=3339042) land on a point value
# that the host's internal points->EMU float32 cast truncates one EMU low
# (e.g. 3896227/12700 -> floored back to 3896226), so those two use a
# tiny (sub-EMU, ~1/100 pt) nudge upward to land on the exact EMU target.
$moves = @{
    "Hexagon 4"  = @{ left = 15.908110236220473; top = 262.4648031496063  }
    "Hexagon 6"  = @{ left = 306.7895355590551;  top = 262.4648031496063  }
    "Hexagon 8"  = @{ left = 597.6709448818898;  top = 13.97236220472441  }
    "Hexagon 9"  = @{ left = 306.7895355590551;  top = 8.078976377952756  }
    "Hexagon 10" = @{ left = 597.6709448818898;  top = 262.9167029133858  }
}

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($moves.ContainsKey($sh.Name)) {
        $m = $moves[$sh.Name]
        $sh.Left = $m.left
        $sh.Top  = $m.top
    }
}
